$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.947566
$ws.Range("H2").Value = 2.842698
$ws.Range("I2").Value = 0.01860557227824198
$ws.Range("J2").Value = 0.01860557227824198
$ws.Range("M2").Value = 45.90594266666667
$ws.Range("N2").Value = 137.717828
$ws.Range("O2").Value = 0.3954672001633582
$ws.Range("P2").Value = 0.3954672001633583
$ws.Range("Q2").Value = 43.49891046888267
$ws.Range("R2").Value = 391.490194219944
$ws.Range("S2").Value = 0.00735789357631335
$ws.Range("T2").Value = 0.007357893576313351

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.947566
$ws.Range("H3").Value = 2.842698
$ws.Range("I3").Value = 0.01860557227824198
$ws.Range("J3").Value = 0.01860557227824198
$ws.Range("O3").Value = 0.3484294080560655
$ws.Range("P3").Value = 0.3484294080560656
$ws.Range("Q3").Value = 38.32504849832267
$ws.Range("R3").Value = 344.925436484904
$ws.Range("S3").Value = 0.006482728535452197
$ws.Range("T3").Value = 0.006482728535452197

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.947566
$ws.Range("H4").Value = 2.842698
$ws.Range("I4").Value = 0.01860557227824198
$ws.Range("J4").Value = 0.01860557227824198
$ws.Range("M4").Value = 12.761795
$ws.Range("N4").Value = 38.28538500000001
$ws.Range("O4").Value = 0.1099393900775594
$ws.Range("P4").Value = 0.1099393900775594
$ws.Range("Q4").Value = 12.09264304097
$ws.Range("R4").Value = 108.83378736873
$ws.Range("S4").Value = 0.002045485268313871
$ws.Range("T4").Value = 0.002045485268313871

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.947566
$ws.Range("H5").Value = 2.842698
$ws.Range("I5").Value = 0.01860557227824198
$ws.Range("J5").Value = 0.01860557227824198
$ws.Range("M5").Value = 16.966758
$ws.Range("N5").Value = 50.900274
$ws.Range("O5").Value = 0.1461640017030168
$ws.Range("P5").Value = 0.1461640017030168
$ws.Range("Q5").Value = 16.077123011028
$ws.Range("R5").Value = 144.694107099252
$ws.Range("S5").Value = 0.002719464898162563
$ws.Range("T5").Value = 0.002719464898162563

# Row 6
$ws.Range("I6").Value = 0.9615746457924133
$ws.Range("J6").Value = 0.9615746457924131
$ws.Range("M6").Value = 45.90594266666667
$ws.Range("N6").Value = 137.717828
$ws.Range("O6").Value = 0.3954672001633582
$ws.Range("P6").Value = 0.3954672001633583
$ws.Range("Q6").Value = 2248.114102643661
$ws.Range("R6").Value = 20233.02692379295
$ws.Range("S6").Value = 0.3802712329195985
$ws.Range("T6").Value = 0.3802712329195986

# Row 7
$ws.Range("I7").Value = 0.9615746457924133
$ws.Range("J7").Value = 0.9615746457924131
$ws.Range("O7").Value = 0.3484294080560655
$ws.Range("P7").Value = 0.3484294080560656
$ws.Range("S7").Value = 0.3350408846351715
$ws.Range("T7").Value = 0.3350408846351715

# Row 8
$ws.Range("I8").Value = 0.9615746457924133
$ws.Range("J8").Value = 0.9615746457924131
$ws.Range("M8").Value = 12.761795
$ws.Range("N8").Value = 38.28538500000001
$ws.Range("O8").Value = 0.1099393900775594
$ws.Range("P8").Value = 0.1099393900775594
$ws.Range("Q8").Value = 624.9729261170318
$ws.Range("R8").Value = 5624.756335053286
$ws.Range("S8").Value = 0.1057149300724631
$ws.Range("T8").Value = 0.1057149300724631

# Row 9
$ws.Range("I9").Value = 0.9615746457924133
$ws.Range("J9").Value = 0.9615746457924131
$ws.Range("M9").Value = 16.966758
$ws.Range("N9").Value = 50.900274
$ws.Range("O9").Value = 0.1461640017030168
$ws.Range("P9").Value = 0.1461640017030168
$ws.Range("Q9").Value = 830.899132448026
$ws.Range("R9").Value = 7478.092192032233
$ws.Range("S9").Value = 0.1405475981651801
$ws.Range("T9").Value = 0.1405475981651801

# Row 10
$ws.Range("G10").Value = 0.8226676666666667
$ws.Range("H10").Value = 2.468003
$ws.Range("I10").Value = 0.01615317849431
$ws.Range("J10").Value = 0.01615317849431
$ws.Range("M10").Value = 45.90594266666667
$ws.Range("N10").Value = 137.717828
$ws.Range("O10").Value = 0.3954672001633582
$ws.Range("P10").Value = 0.3954672001633583
$ws.Range("Q10").Value = 37.76533473972044
$ws.Range("R10").Value = 339.888012657484
$ws.Range("S10").Value = 0.006388052272883745
$ws.Range("T10").Value = 0.006388052272883746

# Row 11
$ws.Range("G11").Value = 0.8226676666666667
$ws.Range("H11").Value = 2.468003
$ws.Range("I11").Value = 0.01615317849431
$ws.Range("J11").Value = 0.01615317849431
$ws.Range("O11").Value = 0.3484294080560655
$ws.Range("P11").Value = 0.3484294080560656
$ws.Range("Q11").Value = 33.27343765289378
$ws.Range("R11").Value = 299.460938876044
$ws.Range("S11").Value = 0.005628242420996401
$ws.Range("T11").Value = 0.005628242420996401

# Row 12
$ws.Range("G12").Value = 0.8226676666666667
$ws.Range("H12").Value = 2.468003
$ws.Range("I12").Value = 0.01615317849431
$ws.Range("J12").Value = 0.01615317849431
$ws.Range("M12").Value = 12.761795
$ws.Range("N12").Value = 38.28538500000001
$ws.Range("O12").Value = 0.1099393900775594
$ws.Range("P12").Value = 0.1099393900775594
$ws.Range("Q12").Value = 10.49871611512834
$ws.Range("R12").Value = 94.48844503615501
$ws.Range("S12").Value = 0.001775870591478391
$ws.Range("T12").Value = 0.001775870591478391

# Row 13
$ws.Range("G13").Value = 0.8226676666666667
$ws.Range("H13").Value = 2.468003
$ws.Range("I13").Value = 0.01615317849431
$ws.Range("J13").Value = 0.01615317849431
$ws.Range("M13").Value = 16.966758
$ws.Range("N13").Value = 50.900274
$ws.Range("O13").Value = 0.1461640017030168
$ws.Range("P13").Value = 0.1461640017030168
$ws.Range("Q13").Value = 13.958003214758
$ws.Range("R13").Value = 125.622028932822
$ws.Range("S13").Value = 0.002361013208951461
$ws.Range("T13").Value = 0.002361013208951461

# Row 14
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 0.186737
$ws.Range("H14").Value = 0.560211
$ws.Range("I14").Value = 0.003666603435034682
$ws.Range("J14").Value = 0.003666603435034681
$ws.Range("M14").Value = 45.90594266666667
$ws.Range("N14").Value = 137.717828
$ws.Range("O14").Value = 0.3954672001633582
$ws.Range("P14").Value = 0.3954672001633583
$ws.Range("Q14").Value = 8.572338015745334
$ws.Range("R14").Value = 77.151042141708
$ws.Range("S14").Value = 0.001450021394562517
$ws.Range("T14").Value = 0.001450021394562517

# Row 15
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 0.186737
$ws.Range("H15").Value = 0.560211
$ws.Range("I15").Value = 0.003666603435034682
$ws.Range("J15").Value = 0.003666603435034681
$ws.Range("O15").Value = 0.3484294080560655
$ws.Range("P15").Value = 0.3484294080560656
$ws.Range("Q15").Value = 7.552724117825334
$ws.Range("R15").Value = 67.974517060428
$ws.Range("S15").Value = 0.001277552464445471
$ws.Range("T15").Value = 0.001277552464445471

# Row 16
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 0.186737
$ws.Range("H16").Value = 0.560211
$ws.Range("I16").Value = 0.003666603435034682
$ws.Range("J16").Value = 0.003666603435034681
$ws.Range("M16").Value = 12.761795
$ws.Range("N16").Value = 38.28538500000001
$ws.Range("O16").Value = 0.1099393900775594
$ws.Range("P16").Value = 0.1099393900775594
$ws.Range("Q16").Value = 2.383099312915
$ws.Range("R16").Value = 21.447893816235
$ws.Range("S16").Value = 0.0004031041453039971
$ws.Range("T16").Value = 0.0004031041453039971

# Row 17
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 0.186737
$ws.Range("H17").Value = 0.560211
$ws.Range("I17").Value = 0.003666603435034682
$ws.Range("J17").Value = 0.003666603435034681
$ws.Range("M17").Value = 16.966758
$ws.Range("N17").Value = 50.900274
$ws.Range("O17").Value = 0.1461640017030168
$ws.Range("P17").Value = 0.1461640017030168
$ws.Range("Q17").Value = 3.168321488646
$ws.Range("R17").Value = 28.514893397814
$ws.Range("S17").Value = 0.0005359254307226965
$ws.Range("T17").Value = 0.0005359254307226965
